$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 15.8494
$ws.Range("A4").Value = -21.50719999999999
$ws.Range("A6").Value = -22.64560000000002
$ws.Range("A7").Value = -19.34899999999999
$ws.Range("C7").Value = -12.06
$ws.Range("A8").Value = -22.37040000000001
$ws.Range("C11").Value = -11.2931
$ws.Range("C12").Value = -10.9554
$ws.Range("D12").Value = -7.442199999999997
$ws.Range("E12").Value = 17.0475
$ws.Range("D13").Value = -8.638599999999999
$ws.Range("E13").Value = 15.8895
$ws.Range("D14").Value = -7.957799999999994
$ws.Range("C15").Value = -14.69779999999999
$ws.Range("A16").Value = -21.60859999999999
$ws.Range("D16").Value = -9.104500000000007
$ws.Range("D19").Value = -8.153600000000001
$ws.Range("A20").Value = -19.8452
$ws.Range("C20").Value = -12.07700000000001
$ws.Range("D20").Value = -7.399000000000004
$ws.Range("A21").Value = -19.86789999999999
$ws.Range("C21").Value = -12.0234
$ws.Range("C22").Value = -12.07199999999999
$ws.Range("D22").Value = -8.426300000000001
$ws.Range("E22").Value = 16.62119999999999
$ws.Range("C23").Value = -11.90720000000001
$ws.Range("E25").Value = 16.99400000000001
$ws.Range("A28").Value = -21.7612
$ws.Range("A29").Value = -21.46629999999998
$ws.Range("C29").Value = -11.64450000000001
$ws.Range("E29").Value = 17.3171
$ws.Range("A30").Value = -21.5717
$ws.Range("A32").Value = -21.2516
$ws.Range("C34").Value = -11.26580000000001
$ws.Range("E34").Value = 17.5963
$ws.Range("D36").Value = -8.703699999999996
$ws.Range("A40").Value = -20.35860000000001
$ws.Range("C42").Value = -11.96739999999999
$ws.Range("C43").Value = -13.66009999999999
$ws.Range("D43").Value = -8.1241
$ws.Range("E43").Value = 16.56929999999999
$ws.Range("C44").Value = -14.0243
$ws.Range("C45").Value = -13.58049999999999
$ws.Range("A46").Value = -22.0303
$ws.Range("C46").Value = -13.0355
$ws.Range("D46").Value = -8.731100000000003
$ws.Range("E48").Value = 17.55640000000001
$ws.Range("C50").Value = -14.10209999999999
$ws.Range("D50").Value = -7.9693
$ws.Range("A51").Value = -21.6251
$ws.Range("C51").Value = -11.00510000000001
$ws.Range("A52").Value = -22.2448
$ws.Range("A57").Value = -22.2061
$ws.Range("C57").Value = -14.17539999999999
$ws.Range("A59").Value = -22.36140000000001
$ws.Range("E60").Value = 15.43290000000001
$ws.Range("A62").Value = -22.28380000000001
$ws.Range("C65").Value = -13.1046
$ws.Range("A66").Value = -21.8132
$ws.Range("C66").Value = -11.2614
$ws.Range("C67").Value = -10.9585
$ws.Range("E68").Value = 17.76490000000002
$ws.Range("E70").Value = 18.23380000000002
$ws.Range("E71").Value = 17.02700000000001
$ws.Range("A73").Value = -20.6273
$ws.Range("E73").Value = 17.44940000000002
$ws.Range("A74").Value = -22.06809999999998
$ws.Range("D76").Value = -7.963400000000006
$ws.Range("A77").Value = -19.71649999999999
$ws.Range("E78").Value = 17.16900000000002
$ws.Range("C79").Value = -12.13930000000001
$ws.Range("C84").Value = -13.26369999999999
$ws.Range("C87").Value = -13.4877
$ws.Range("E87").Value = 16.16709999999998
$ws.Range("A92").Value = -21.5507
$ws.Range("C92").Value = -11.27660000000001
$ws.Range("E92").Value = 18.29360000000002
$ws.Range("D95").Value = -8.161200000000001
$ws.Range("C97").Value = -11.5786
$ws.Range("D97").Value = -8.766299999999996
$ws.Range("D99").Value = -8.289600000000005
$ws.Range("A100").Value = -22.2042
$ws.Range("E101").Value = 17.14420000000001
